$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task rows (85-87), formerly blank placeholder rows ---

# Row 85
$ws.Range("A85").Value2 = 84
$ws.Range("B85").Value2 = "Follow up with Sarah Elkins on the status of updating AHP on caArray STAGE and TRAINING to use GitHub."
$ws.Range("C85").Value2 = "Mike Hunter"
$ws.Range("D85").Value2 = 39946
$ws.Range("E85").Value2 = "Assigned"

# Row 86
$ws.Range("A86").Value2 = 85
$ws.Range("B86").Value2 = "Contact Eva Shalley to find out if TRANSCEND plans to upgrade to the next release of caIntegrator."
$ws.Range("C86").Value2 = "Mike Hunter"
$ws.Range("D86").Value2 = 39946
$ws.Range("E86").Value2 = "Assigned"

# Row 87
$ws.Range("A87").Value2 = 86
$ws.Range("B87").Value2 = "Invite Laxmi Lolla to attend an upcoming status meeting to share her feedback and recommendations regarding usability and performance."
$ws.Range("C87").Value2 = "Ulrike Wagner"
$ws.Range("D87").Value2 = 39946
$ws.Range("E87").Value2 = "Assigned"

# Rows 85-87 grew taller (wrapped, 2-line content) in real Excel; best-effort
# reproduction of that auto row height via an explicit height.
$ws.Range("A85:E87").RowHeight = 31

# --- Renumber the remaining blank placeholder rows (88-92), column A only ---
$ws.Range("A88").Value2 = 87
$ws.Range("A89").Value2 = 88
$ws.Range("A90").Value2 = 89
$ws.Range("A91").Value2 = 90
$ws.Range("A92").Value2 = 91

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("B88").Select()

# --- Shift the workbook window horizontally on screen ---
$win = $excel.Windows.Item(1)
$win.Left = 280
